$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 1.211767666666667
$ws.Range("H2").Value = 3.635303
$ws.Range("I2").Value = 0.0191872416143265
$ws.Range("J2").Value = 0.0191872416143265
$ws.Range("M2").Value = 0.3620403333333334
$ws.Range("N2").Value = 1.086121
$ws.Range("O2").Value = 0.138293228945037
$ws.Range("P2").Value = 0.138293228945037
$ws.Range("Q2").Value = 0.4387087699625555
$ws.Range("R2").Value = 3.948378929663
$ws.Range("S2").Value = 0.002653465597393796
$ws.Range("T2").Value = 0.002653465597393796
$ws.Range("G3").Value = 1.211767666666667
$ws.Range("H3").Value = 3.635303
$ws.Range("I3").Value = 0.0191872416143265
$ws.Range("J3").Value = 0.0191872416143265
$ws.Range("O3").Value = 0.530310855165568
$ws.Range("P3").Value = 0.530310855165568
$ws.Range("Q3").Value = 1.682309573232555
$ws.Range("R3").Value = 15.140786159093
$ws.Range("S3").Value = 0.01017520250876186
$ws.Range("T3").Value = 0.01017520250876186
$ws.Range("G4").Value = 1.211767666666667
$ws.Range("H4").Value = 3.635303
$ws.Range("I4").Value = 0.0191872416143265
$ws.Range("J4").Value = 0.0191872416143265
$ws.Range("M4").Value = 0.8675673333333332
$ws.Range("N4").Value = 2.602702
$ws.Range("O4").Value = 0.331395915889395
$ws.Range("P4").Value = 0.331395915889395
$ws.Range("Q4").Value = 1.051290043189555
$ws.Range("R4").Value = 9.461610388705997
$ws.Range("S4").Value = 0.006358573508170845
$ws.Range("T4").Value = 0.006358573508170845
$ws.Range("G5").Value = 27.75404733333334
$ws.Range("H5").Value = 83.26214200000001
$ws.Range("I5").Value = 0.4394601594090953
$ws.Range("J5").Value = 0.4394601594090954
$ws.Range("M5").Value = 0.3620403333333334
$ws.Range("N5").Value = 1.086121
$ws.Range("O5").Value = 0.138293228945037
$ws.Range("P5").Value = 0.138293228945037
$ws.Range("Q5").Value = 10.04808454790911
$ws.Range("R5").Value = 90.43276093118202
$ws.Range("S5").Value = 0.06077436443738447
$ws.Range("T5").Value = 0.06077436443738448
$ws.Range("G6").Value = 27.75404733333334
$ws.Range("H6").Value = 83.26214200000001
$ws.Range("I6").Value = 0.4394601594090953
$ws.Range("J6").Value = 0.4394601594090954
$ws.Range("O6").Value = 0.530310855165568
$ws.Range("P6").Value = 0.530310855165568
$ws.Range("R6").Value = 346.7810763422021
$ws.Range("S6").Value = 0.2330504929474342
$ws.Range("T6").Value = 0.2330504929474342
$ws.Range("G7").Value = 27.75404733333334
$ws.Range("H7").Value = 83.26214200000001
$ws.Range("I7").Value = 0.4394601594090953
$ws.Range("J7").Value = 0.4394601594090954
$ws.Range("M7").Value = 0.8675673333333332
$ws.Range("N7").Value = 2.602702
$ws.Range("O7").Value = 0.331395915889395
$ws.Range("P7").Value = 0.331395915889395
$ws.Range("Q7").Value = 24.07850483418711
$ws.Range("R7").Value = 216.706543507684
$ws.Range("S7").Value = 0.1456353020242767
$ws.Range("T7").Value = 0.1456353020242767
$ws.Range("G8").Value = 34.18905066666667
$ws.Range("H8").Value = 102.567152
$ws.Range("I8").Value = 0.5413525989765782
$ws.Range("J8").Value = 0.5413525989765782
$ws.Range("M8").Value = 0.3620403333333334
$ws.Range("N8").Value = 1.086121
$ws.Range("O8").Value = 0.138293228945037
$ws.Range("P8").Value = 0.138293228945037
$ws.Range("Q8").Value = 12.37781529971022
$ws.Range("R8").Value = 111.400337697392
$ws.Range("S8").Value = 0.07486539891025872
$ws.Range("T8").Value = 0.07486539891025872
$ws.Range("G9").Value = 34.18905066666667
$ws.Range("H9").Value = 102.567152
$ws.Range("I9").Value = 0.5413525989765782
$ws.Range("J9").Value = 0.5413525989765782
$ws.Range("O9").Value = 0.530310855165568
$ws.Range("P9").Value = 0.530310855165568
$ws.Range("Q9").Value = 47.46501232739022
$ws.Range("R9").Value = 427.185110946512
$ws.Range("S9").Value = 0.2870851597093719
$ws.Range("T9").Value = 0.2870851597093719
$ws.Range("G10").Value = 34.18905066666667
$ws.Range("H10").Value = 102.567152
$ws.Range("I10").Value = 0.5413525989765782
$ws.Range("J10").Value = 0.5413525989765782
$ws.Range("M10").Value = 0.8675673333333332
$ws.Range("N10").Value = 2.602702
$ws.Range("O10").Value = 0.331395915889395
$ws.Range("P10").Value = 0.331395915889395
$ws.Range("Q10").Value = 29.66130351607822
$ws.Range("R10").Value = 266.951731644704
$ws.Range("S10").Value = 0.1794020403569475
$ws.Range("T10").Value = 0.1794020403569475
